$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date serial numbers (Excel 1900 date system) corresponding to the
# former "YYYYQ4" text labels in column A (Dec 31 of each year).
$dates = @(
    32142, 32508, 32873, 33238, 33603, 33969, 34334, 34699, 35064, 35430,
    35795, 36160, 36525, 36891, 37256, 37621, 37986, 38352, 38717, 39082,
    39447, 39813, 40178, 40543, 40908, 41274, 41639, 42004, 42369, 42735,
    43100, 43465, 43830, 44196, 44561, 44926, 45291, 45657
)

$row = 2
foreach ($d in $dates) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $d
    $cell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $row = $row + 1
}
